# Update gh-pages output data: refresh "想去人数" (want-to-go count) figures
# and one cover image URL across the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F3").Value = 3246
    $ws.Range("F9").Value = 471
    $ws.Range("F10").Value = 376
    $ws.Range("F21").Value = 58
    $ws.Range("F24").Value = 229
    $ws.Range("F26").Value = 39
    $ws.Range("F29").Value = 325
    $ws.Range("F30").Value = 2222
    $ws.Range("F34").Value = 443
    $ws.Range("F38").Value = 349
    $ws.Range("F40").Value = 525

    $ws.Range("I37").Value = "//i2.hdslb.com/bfs/openplatform/202406/65hJjOfJ1717642614493.jpeg"
}
